# Generate Report for Handback
# For each locale sheet (zh-cn, de-de), populate the "Latest Target File" (F)
# and "Latest Handback File" (G) columns with the same file references as the
# "Latest Handoff File" column, mark the status as handed back, and stamp the
# handback datetime.

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"
$localeSheets = @("zh-cn", "de-de")

# The "Overview" sheet mirrors each locale's status in columns B (zh-cn) and
# C (de-de); the same "Ready for handoff" -> "Handed back: in sync with
# en-US" text change applies there for both rows.
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = $statusText
$overview.Range("C2").Value = $statusText
$overview.Range("B3").Value = $statusText
$overview.Range("C3").Value = $statusText

foreach ($sheetName in $localeSheets) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Snapshot the existing hyperlinks (source markdown + handoff xlf) for
    # row 2 before we start adding new ones, since the collection mutates.
    $existing = @($ws.Hyperlinks)
    $mdLink = $null
    $xlfLink = $null
    foreach ($h in $existing) {
        $refAddr = $h.Range.Address(0, 0)
        if ($refAddr -eq "F2" -or $refAddr -eq "1048576") {}
    }

    foreach ($rowNum in 2, 3) {
        $aAddr = "A" + $rowNum
        $dAddr = "D" + $rowNum
        $fAddr = "F" + $rowNum
        $gAddr = "G" + $rowNum
        $hAddr = "H" + $rowNum
        $cAddr = "C" + $rowNum

        $mdUrl = ""
        $mdDisplay = ""
        $xlfUrl = ""
        $xlfDisplay = ""
        foreach ($h in @($ws.Hyperlinks)) {
            if ($h.Range.Row -eq 2 -and $h.Range.Column -eq 1) {
                $mdUrl = $h.Address
                $mdDisplay = $h.TextToDisplay
            }
            if ($h.Range.Row -eq 2 -and $h.Range.Column -eq 4) {
                $xlfUrl = $h.Address
                $xlfDisplay = $h.TextToDisplay
            }
        }

        # Status -> handed back
        $ws.Range($cAddr).Value = $statusText

        # Latest Target File (F) = same markdown source as column A
        $ws.Range($fAddr).Value = $mdDisplay
        $ws.Hyperlinks.Add($ws.Range($fAddr), $mdUrl, "", "", $mdDisplay)
        $ws.Range($fAddr).Style = $ws.Range("A" + $rowNum).Style

        # Latest Handback File (G) = same xlf file as column D
        $ws.Range($gAddr).Value = $xlfDisplay
        $ws.Hyperlinks.Add($ws.Range($gAddr), $xlfUrl, "", "", $xlfDisplay)
        $ws.Range($gAddr).Style = $ws.Range("D" + $rowNum).Style
    }

    # Latest Handback DateTime (H) stamped with the handback timestamp.
    if ($sheetName -eq "zh-cn") {
        $ws.Range("H2").Value = "2016-03-22 15:11:26"
        $ws.Range("H3").Value = "2016-03-22 15:11:26"
    } else {
        $ws.Range("H2").Value = "2016-03-22 15:11:36"
        $ws.Range("H3").Value = "2016-03-22 15:11:36"
    }
}

Write-Host "done"
